$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $cell = $ws.Range($rangeAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "53.928.46"
Set-TextValue "E2" "  -1.12%  "
Set-TextValue "D3" "2.253.87"
Set-TextValue "E3" "  +0.38%  "
Set-TextValue "E4" "  +0.82%  "
Set-TextValue "D5" "493.62"
Set-TextValue "E5" "  -0.53%  "
Set-TextValue "D6" "128.07"
Set-TextValue "E6" "  +0.28%  "
Set-TextValue "E7" "  +0.50%  "
Set-TextValue "E8" "  -1.30%  "
Set-TextValue "E9" "  -0.44%  "
Set-TextValue "E10" "  +0.93%  "
Set-TextValue "E11" "  +2.83%  "
Set-TextValue "E12" "  +1.93%  "
Set-TextValue "D13" "2.653.49"
Set-TextValue "E13" "  -0.49%  "
Set-TextValue "D14" "22.62"
Set-TextValue "E14" "  +3.75%  "
Set-TextValue "D15" "53.905.95"
Set-TextValue "E15" "  -0.64%  "
Set-TextValue "E16" "  -0.24%  "
Set-TextValue "D17" "2.252.95"
Set-TextValue "E17" "  -1.25%  "
Set-TextValue "D18" "10.21"
Set-TextValue "E18" "  +1.45%  "
Set-TextValue "E19" "  +0.14%  "
Set-TextValue "D20" "301.04"
Set-TextValue "E20" "  -0.36%  "
Set-TextValue "E21" "  -3.00%  "
Set-TextValue "E22" "  +0.62%  "
Set-TextValue "D23" "60.62"
Set-TextValue "E23" "  -3.04%  "
Set-TextValue "D24" "0.998"
Set-TextValue "E24" "  +0.04%  "
Set-TextValue "E25" "  -1.79%  "
Set-TextValue "E26" "  +2.37%  "
Set-TextValue "D27" "170.23"
Set-TextValue "E27" "  +0.52%  "
Set-TextValue "E28" "  -0.55%  "
Set-TextValue "E29" "  +0.38%  "
Set-TextValue "D30" "0.0₃0685"
Set-TextValue "E30" "  -0.62%  "
Set-TextValue "E31" "  +0.25%  "
Set-TextValue "E32" "  +0.27%  "
Set-TextValue "D33" "17.75"
Set-TextValue "E33" "  +0.77%  "
Set-TextValue "D34" "0.999"
Set-TextValue "E34" "  +0.71%  "
Set-TextValue "D35" "0.938"
Set-TextValue "E35" "  +8.06%  "
Set-TextValue "E36" "  -0.69%  "
Set-TextValue "E37" "  -1.05%  "
Set-TextValue "E38" "  -1.48%  "
Set-TextValue "E39" "  -1.61%  "
Set-TextValue "E40" "  -0.12%  "
Set-TextValue "D41" "125.54"
Set-TextValue "E41" "  -2.26%  "
Set-TextValue "D42" "4.76"
Set-TextValue "E42" "  -3.57%  "
Set-TextValue "E43" "  +0.65%  "
Set-TextValue "E44" "  -0.46%  "
Set-TextValue "D45" "0.541"
Set-TextValue "E45" "  -0.52%  "
Set-TextValue "D46" "238.76"
Set-TextValue "E46" "  -0.62%  "
Set-TextValue "D47" "0.368"
Set-TextValue "E47" "  -1.43%  "
Set-TextValue "E48" "  -0.25%  "
Set-TextValue "E49" "  +0.31%  "
Set-TextValue "D50" "16.08"
Set-TextValue "E50" "  -1.77%  "
Set-TextValue "E51" "  -1.05%  "
